$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7 (pushes old row 7.. down by one).
# The new row inherits formatting from the row above (row 6), matching
# Excel's default "Format Same As Above" behavior for Insert Sheet Rows.
$ws.Rows.Item(7).Insert()

# Fill the new row 7 with the new "Apochromat 40x" objective-lens entry
# (same Camera/Microscope/Magnification/Width/Height/Scale/Unit as the
# neighboring Neofluar 40x row -- only the objective name and the pixel
# distance differ).
$ws.Range("A7").Value = "Canon EOS RP"
$ws.Range("B7").Value = "Zeiss Axio Scope A1"
$ws.Range("C7").Value = "Apochromat 40x"
$ws.Range("D7").Value = "40x"
$ws.Range("E7").Value = 6240
$ws.Range("F7").Value = 4160
$ws.Range("G7").Value = "873"
$ws.Range("H7").Value = 50
$ws.Range("I7").Value = "µm"

# Row 6 (Neofluar 40x) is no longer the last entry of the "Zeiss Axio
# Scope A1" group, so its bottom border is cleared (the new row 7 now
# carries the group-end divider that Insert copied down from row 6).
$ws.Range("A6:I6").Borders.LineStyle = -4142

$ws.Range("E6").Select()
